$wb = $excel.ActiveWorkbook

# --- Sam sheet: fill in time card entries for week 1 (rows 2-4) ---
$samWs = $wb.Worksheets.Item("Sam")

$samWs.Range("A2").Value = 42979.458333333336
$samWs.Range("B2").Value = 42979.5

$samWs.Range("A3").Value = 42979.354166666664
$samWs.Range("B3").Value = 42979.4375

$samWs.Range("A4").Value = 42982.458333333336
$samWs.Range("B4").Value = 42982.583333333336

# Make "Sam" the active sheet/tab and leave the selection where the author left it
$samWs.Activate()
[void]$samWs.Range("E17").Select()

# --- Travis sheet: move the selection (no data changes here) ---
$travisWs = $wb.Worksheets.Item("Travis")
[void]$travisWs.Range("C8").Select()

# Re-activate Sam so it is the sheet shown/selected when the workbook is saved
$samWs.Activate()
